# Auto-generated Excel COM-interop script applying Goblin_Profits scheduled-runner updates
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 262.62857
$ws.Range("I53").Value = 181.42857
$ws.Range("J53").Value = 587.4286
$ws.Range("K53").Value = 181.42857
$ws.Range("L53").Value = 587.4286
$ws.Range("M53").Value = 455.57143
$ws.Range("N53").Value = -1861.4286

# Row 94
$ws.Range("H94").Value = 2966
$ws.Range("I94").Value = 3095.5557
$ws.Range("J94").Value = 1800
$ws.Range("K94").Value = 3095.5557
$ws.Range("L94").Value = 1800
$ws.Range("M94").Value = -2644.5557
$ws.Range("N94").Value = -2702

# Row 101
$ws.Range("H101").Value = 2549.3076
$ws.Range("J101").Value = 14999.5
$ws.Range("L101").Value = 44998.5
$ws.Range("N101").Value = -48242.5

# Row 105
$ws.Range("H105").Value = 72500
$ws.Range("J105").Value = 72500
$ws.Range("L105").Value = 72500
$ws.Range("N105").Value = -79488

# Row 113
$ws.Range("H113").Value = 3699.2
$ws.Range("J113").Value = 4856.4287
$ws.Range("L113").Value = 4856.4287
$ws.Range("N113").Value = -11364.4287

# Row 121
$ws.Range("H121").Value = 2895.7778
$ws.Range("J121").Value = 2895.7778
$ws.Range("L121").Value = 8687.3334
$ws.Range("N121").Value = -12181.3334

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1687.5
$ws.Range("I2").Value = 654.5
$ws.Range("J2").Value = 3340.3
$ws.Range("K2").Value = 654.5
$ws.Range("L2").Value = 3340.3
$ws.Range("M2").Value = -541.5
$ws.Range("N2").Value = -3566.3

# Row 45
$ws.Range("H45").Value = 2666.5
$ws.Range("I45").Value = 1083.6
$ws.Range("K45").Value = 1083.6
$ws.Range("M45").Value = -706.5999999999999

# Row 116
$ws.Range("H116").Value = 1687.5
$ws.Range("I116").Value = 654.5
$ws.Range("J116").Value = 3340.3
$ws.Range("K116").Value = 654.5
$ws.Range("L116").Value = 3340.3
$ws.Range("M116").Value = 1639.5
$ws.Range("N116").Value = -7928.3

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1687.5
$ws.Range("I3").Value = 654.5
$ws.Range("J3").Value = 3340.3
$ws.Range("K3").Value = 654.5
$ws.Range("L3").Value = 3340.3
$ws.Range("M3").Value = -540.5
$ws.Range("N3").Value = -3568.3

# Row 53
$ws.Range("H53").Value = 99833.336
$ws.Range("I53").Value = 99500
$ws.Range("K53").Value = 99500
$ws.Range("M53").Value = -98926

# Row 105
$ws.Range("H105").Value = 2393.5715
$ws.Range("I105").Value = 1438.75
$ws.Range("K105").Value = 1438.75
$ws.Range("M105").Value = 308.25

$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 69995
$ws.Range("J50").Value = 69995
$ws.Range("L50").Value = 69995
$ws.Range("N50").Value = -71245

# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# Row 58
$ws.Range("H58").Value = 2044.3334
$ws.Range("I58").Value = 2063
$ws.Range("K58").Value = 2063
$ws.Range("M58").Value = -1860

# Row 59
$ws.Range("H59").Value = 38728.75
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 38728.75
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 38728.75
$ws.Range("N59").Value = -41018.75
$ws.Range("M59").ClearContents()

# Row 60
$ws.Range("H60").Value = 18400
$ws.Range("I60").Value = 13000
$ws.Range("J60").Value = 40000
$ws.Range("K60").Value = 13000
$ws.Range("L60").Value = 40000
$ws.Range("M60").Value = -12489
$ws.Range("N60").Value = -41022

# Row 61
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# Row 68
$ws.Range("H68").Value = 64999
$ws.Range("I68").Value = 59998.5
$ws.Range("K68").Value = 59998.5
$ws.Range("M68").Value = -59249.5

# Row 71
$ws.Range("H71").Value = 64999
$ws.Range("I71").Value = 59998.5
$ws.Range("K71").Value = 179995.5
$ws.Range("M71").Value = -176251.5

# Row 136
$ws.Range("H136").Value = 2044.3334
$ws.Range("I136").Value = 2063
$ws.Range("K136").Value = 6189
$ws.Range("M136").Value = -3639

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 2589.5483
$ws.Range("I68").Value = 2389.4443
$ws.Range("J68").Value = 2671.4092
$ws.Range("K68").Value = 7168.3329
$ws.Range("L68").Value = 8014.2276
$ws.Range("M68").Value = -6357.3329
$ws.Range("N68").Value = -9636.2276

# Row 71
$ws.Range("H71").Value = 2589.5483
$ws.Range("I71").Value = 2389.4443
$ws.Range("J71").Value = 2671.4092
$ws.Range("K71").Value = 21504.9987
$ws.Range("L71").Value = 24042.6828
$ws.Range("M71").Value = -17448.9987
$ws.Range("N71").Value = -32154.6828

# Row 80
$ws.Range("H80").Value = 2437.6667
$ws.Range("I80").Value = 2407.6667
$ws.Range("J80").Value = 2467.6667
$ws.Range("K80").Value = 7223.000100000001
$ws.Range("L80").Value = 7403.000100000001
$ws.Range("M80").Value = -6287.000100000001
$ws.Range("N80").Value = -9275.000100000001

# Row 83
$ws.Range("H83").Value = 2437.6667
$ws.Range("I83").Value = 2407.6667
$ws.Range("J83").Value = 2467.6667
$ws.Range("K83").Value = 21669.0003
$ws.Range("L83").Value = 22209.0003
$ws.Range("M83").Value = -16989.0003
$ws.Range("N83").Value = -31569.0003

# Row 92
$ws.Range("H92").Value = 1634.6666
$ws.Range("J92").Value = 866.3333
$ws.Range("L92").Value = 2598.9999
$ws.Range("N92").Value = -5094.9999

# Row 129
$ws.Range("H129").Value = 2956.2354
$ws.Range("J129").Value = 4297.3335
$ws.Range("L129").Value = 12892.0005
$ws.Range("N129").Value = -22892.0005

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 19856.715
$ws.Range("I80").Value = 29500.5
$ws.Range("J80").Value = 6998.3335
$ws.Range("K80").Value = 29500.5
$ws.Range("L80").Value = 6998.3335
$ws.Range("M80").Value = -28502.5
$ws.Range("N80").Value = -8994.333500000001

# Row 83
$ws.Range("H83").Value = 19856.715
$ws.Range("I83").Value = 29500.5
$ws.Range("J83").Value = 6998.3335
$ws.Range("K83").Value = 147502.5
$ws.Range("L83").Value = 34991.6675
$ws.Range("M83").Value = -142510.5
$ws.Range("N83").Value = -44975.6675

# Row 97
$ws.Range("H97").Value = 8564.058999999999
$ws.Range("I97").Value = 2723.2727
$ws.Range("J97").Value = 19272.166
$ws.Range("K97").Value = 2723.2727
$ws.Range("L97").Value = 19272.166
$ws.Range("M97").Value = -2227.2727
$ws.Range("N97").Value = -20264.166

# Row 122
$ws.Range("H122").Value = 6911.391
$ws.Range("I122").Value = 7742.8823
$ws.Range("K122").Value = 23228.6469
$ws.Range("M122").Value = -20778.6469

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2932.4062
$ws.Range("J22").Value = 3269
$ws.Range("L22").Value = 3269
$ws.Range("N22").Value = -3859

# Row 27
$ws.Range("H27").Value = 2932.4062
$ws.Range("J27").Value = 3269
$ws.Range("L27").Value = 3269
$ws.Range("N27").Value = -3483

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 5223.143
$ws.Range("I122").Value = 2898.923
$ws.Range("K122").Value = 8696.769
$ws.Range("M122").Value = -6246.769

# Row 137
$ws.Range("H137").Value = 70715
$ws.Range("J137").Value = 70715
$ws.Range("L137").Value = 70715
$ws.Range("N137").Value = -80915
